$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (id=65) describing the new "Conversion efficiency" parameter.
$newRow = 66
$ws.Cells.Item($newRow, 1).Value = 65
$ws.Cells.Item($newRow, 2).Value = "Conversion efficiency"
$ws.Cells.Item($newRow, 3).Value = "Conversion efficiency factors for electricity and heat power plants, as well as hydrogen production"

# Mirror the author's navigation: scroll down and land the selection near the new row.
$ws.Range("A55").Select()
$ws.Range("C68").Select()
